$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text prefix ("'") ensures values like "24.42" or "0.9978" stay
# text cells (matching the source inlineStr type) instead of being
# auto-converted to numbers by Excel's value parser.

$ws.Range("D2").Value = "'29.368.46"
$ws.Range("E2").Value = "'  -0.14%  "
$ws.Range("D3").Value = "'1.844.22"
$ws.Range("D4").Value = "'0.9978"
$ws.Range("E4").Value = "'  -0.35%  "
$ws.Range("D5").Value = "'240.35"
$ws.Range("D6").Value = "'0.6318"
$ws.Range("E6").Value = "'  +0.71%  "
$ws.Range("D7").Value = "'0.9988"
$ws.Range("E7").Value = "'  -0.27%  "
$ws.Range("E8").Value = "'  -1.78%  "
$ws.Range("E9").Value = "'  +0.00%  "
$ws.Range("D10").Value = "'24.42"
$ws.Range("E10").Value = "'  -1.35%  "
$ws.Range("D11").Value = "'0.07711"
$ws.Range("E11").Value = "'  -0.49%  "
$ws.Range("D12").Value = "'1.844.07"
$ws.Range("E12").Value = "'  -2.34%  "
$ws.Range("E13").Value = "'  -0.63%  "
$ws.Range("D14").Value = "'0.6801"
$ws.Range("E14").Value = "'  +0.17%  "
$ws.Range("D15").Value = "'0.00001028"
$ws.Range("E15").Value = "'  -3.57%  "
$ws.Range("D16").Value = "'82.17"
$ws.Range("E16").Value = "'  -1.29%  "
$ws.Range("D17").Value = "'2.105.08"
$ws.Range("E17").Value = "'  -3.74%  "
$ws.Range("D18").Value = "'6.156"
$ws.Range("E18").Value = "'  -0.04%  "
$ws.Range("D19").Value = "'29.376.56"
$ws.Range("E19").Value = "'  -0.18%  "
$ws.Range("D20").Value = "'229.30"
$ws.Range("E20").Value = "'  +1.27%  "
$ws.Range("E21").Value = "'  +0.07%  "
$ws.Range("D22").Value = "'0.9988"
$ws.Range("E22").Value = "'  -0.28%  "
$ws.Range("D23").Value = "'7.460"
$ws.Range("E23").Value = "'  -0.18%  "
$ws.Range("D24").Value = "'0.9989"
$ws.Range("E24").Value = "'  -0.31%  "
$ws.Range("D25").Value = "'158.83"
$ws.Range("E25").Value = "'  +0.57%  "
$ws.Range("E26").Value = "'  +0.01%  "
$ws.Range("D27").Value = "'8.414"
$ws.Range("E27").Value = "'  -0.05%  "
$ws.Range("E28").Value = "'  -0.69%  "
$ws.Range("D29").Value = "'0.06367"
$ws.Range("E29").Value = "'  +13.97%  "
$ws.Range("E30").Value = "'  -0.11%  "
$ws.Range("D31").Value = "'1.473"
$ws.Range("E31").Value = "'  +0.40%  "
$ws.Range("E32").Value = "'  -0.82%  "
$ws.Range("D33").Value = "'4.062"
$ws.Range("E33").Value = "'  +0.14%  "
$ws.Range("E34").Value = "'  -0.93%  "
$ws.Range("E35").Value = "'  -1.88%  "
$ws.Range("D36").Value = "'0.7002"
$ws.Range("E36").Value = "'  +0.81%  "
$ws.Range("E37").Value = "'  -0.47%  "
$ws.Range("D38").Value = "'2.834"
$ws.Range("E38").Value = "'  +4.18%  "
$ws.Range("D39").Value = "'1.256.31"
$ws.Range("E39").Value = "'  +1.88%  "
$ws.Range("D40").Value = "'0.01823"
$ws.Range("E40").Value = "'  +1.18%  "
$ws.Range("D41").Value = "'6.598"
$ws.Range("E41").Value = "'  +2.88%  "
$ws.Range("D42").Value = "'0.9059"
$ws.Range("E42").Value = "'  +0.11%  "
$ws.Range("D43").Value = "'0.9981"
$ws.Range("D44").Value = "'2.005.77"
$ws.Range("E44").Value = "'  -18.49%  "
$ws.Range("E45").Value = "'  -0.30%  "
$ws.Range("D46").Value = "'66.36"
$ws.Range("E46").Value = "'  +0.61%  "
$ws.Range("E47").Value = "'  -1.99%  "
$ws.Range("D48").Value = "'0.1181"
$ws.Range("D49").Value = "'7.050"
$ws.Range("E49").Value = "'  -1.72%  "
$ws.Range("D50").Value = "'1.704"
$ws.Range("E50").Value = "'  +1.52%  "
$ws.Range("D51").Value = "'9.046"
$ws.Range("E51").Value = "'  +0.71%  "
